# logo_ttalkr.pptx - "updated ttlight and tested new ttgrowth (v2)"
#
# On slide 1, the existing group "Gruppieren 7" (the hexagon + picture +
# "ttalkR" textbox logo) and the standalone "Oval 8" shape are selected
# together and grouped into a brand-new outer group ("Gruppieren 1"),
# exactly like using Format > Group > Group in the UI. The old group keeps
# its id/name/contents untouched and simply becomes a nested child of the
# new wrapper group.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$range = $s.Shapes.Range(@("Gruppieren 7", "Oval 8"))
$newGroup = $range.Group()
$newGroup.Name = "Gruppieren 1"
